$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" footer timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Junio de 2020 a las 01:37"

# Country name re-rank swaps (column A) caused by refreshed case-count sort order
$ws.Range("A80").Value = "Guinea"
$ws.Range("A81").Value = "Haiti"
$ws.Range("A84").Value = "Gabon"
$ws.Range("A85").Value = "Kenia"
$ws.Range("A151").Value = "Libia"
$ws.Range("A152").Value = "Tanzania"
$ws.Range("A153").Value = "Reunion"
$ws.Range("A206").Value = "Islas Malvinas"
$ws.Range("A207").Value = "Groenlandia"
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("A210").Value = "Seychelles"
$ws.Range("A211").Value = "Montserrat"
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("A214").Value = "Islas Virgenes Britanicas"

# Refreshed numeric stats (columns B-H)
$ws.Range("B4").Value = 2261978
$ws.Range("C4").Value = 26251
$ws.Range("D4").Value = 929752
$ws.Range("E4").Value = 1211603
$ws.Range("G4").Value = 682
$ws.Range("H4").Value = 120623
$ws.Range("B5").Value = 983359
$ws.Range("C5").Value = 23050
$ws.Range("E5").Value = 431983
$ws.Range("G5").Value = 1204
$ws.Range("H5").Value = 47869
$ws.Range("B14").Value = 190126
$ws.Range("C14").Value = 622
$ws.Range("E14").Value = 7080
$ws.Range("G14").Value = 19
$ws.Range("H14").Value = 8946
$ws.Range("B21").Value = 100220
$ws.Range("C21").Value = 367
$ws.Range("D21").Value = 62496
$ws.Range("E21").Value = 29424
$ws.Range("G21").Value = 46
$ws.Range("H21").Value = 8300
$ws.Range("B26").Value = 60217
$ws.Range("C26").Value = 3171
$ws.Range("D26").Value = 22680
$ws.Range("E26").Value = 35587
$ws.Range("G26").Value = 86
$ws.Range("H26").Value = 1950
$ws.Range("B53").Value = 18480
$ws.Range("C53").Value = 745
$ws.Range("D53").Value = 6307
$ws.Range("E53").Value = 11698
$ws.Range("G53").Value = 6
$ws.Range("H53").Value = 475
$ws.Range("B66").Value = 10280
$ws.Range("C66").Value = 118
$ws.Range("D66").Value = 7440
$ws.Range("E66").Value = 2506
$ws.Range("B69").Value = 8708
$ws.Range("C69").Value = 16
$ws.Range("E69").Value = 326
$ws.Range("B80").Value = 4841
$ws.Range("C80").Value = 173
$ws.Range("D80").Value = 3467
$ws.Range("E80").Value = 1348
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 26
$ws.Range("B81").Value = 4688
$ws.Range("C81").Value = 141
$ws.Range("D81").Value = 24
$ws.Range("E81").Value = 4582
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 82
$ws.Range("B84").Value = 4340
$ws.Range("C84").Value = 111
$ws.Range("D84").Value = 1657
$ws.Range("E84").Value = 2651
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 32
$ws.Range("B85").Value = 4257
$ws.Range("C85").Value = 213
$ws.Range("D85").Value = 1459
$ws.Range("E85").Value = 2681
$ws.Range("G85").Value = 10
$ws.Range("H85").Value = 117
$ws.Range("D126").Value = 901
$ws.Range("E126").Value = 52
$ws.Range("B135").Value = 850
$ws.Range("C135").Value = 1
$ws.Range("D135").Value = 814
$ws.Range("E135").Value = 12
$ws.Range("B149").Value = 547
$ws.Range("C149").Value = 3
$ws.Range("E149").Value = 181
$ws.Range("B151").Value = 510
$ws.Range("C151").Value = 10
$ws.Range("D151").Value = 81
$ws.Range("E151").Value = 419
$ws.Range("H151").Value = 10
$ws.Range("B152").Value = 509
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 183
$ws.Range("E152").Value = 305
$ws.Range("H152").Value = 21
$ws.Range("B153").Value = 502
$ws.Range("C153").Value = 5
$ws.Range("D153").Value = 460
$ws.Range("E153").Value = 41
$ws.Range("H153").Value = 1
$ws.Range("D178").Value = 74
$ws.Range("E178").Value = 19
$ws.Range("D193").Value = 26
$ws.Range("E193").Value = 3
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
